# Daily attendance processing - 2025-10-20 23:18:15
# Normalize "Recorded By" (column G) entries: when a recorder pair is
# stored as "<user>, System" or "<user>, admin@admin.com", flip it so the
# automated recorder name comes first, e.g. "System, <user>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "

        if ($parts.Count -eq 2 -and $parts[0] -ne "System") {
            if ($parts[1] -eq "System" -or $parts[1] -eq "admin@admin.com") {
                $cell.Value2 = $parts[1] + ", " + $parts[0]
            }
        }
    }
}
